# Rename the worksheet "cond_eb1_c" -> "cond"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "cond"

# Work-in-progress setup of iterative calculation (per commit message:
# "Started implementing [eb2]. Not finished.") — the author updated the
# "Maximum Change" value used for iterative calculation in
# File > Options > Formulas without (yet) enabling iterative calculation
# itself, which is why only the MaxChange/iterateDelta value changes.
$excel.MaxChange = 0.0001
